$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell E1 - set value then copy formatting from an existing header cell (D1)
$ws.Range("E1").Value = "historico"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# Update existing row 2 values
$ws.Range("A2").Value = "a "
$ws.Range("B2").Value = "a "
$ws.Range("C2").Value = "0:00:09.036521"
$ws.Range("D2").Value = "Finalizada"
$ws.Range("E2").Value = "Iniciada em: 23:10:01`nPausada em: 23:10:06`nIniciada em: 23:14:57`nPausada em: 23:15:02`nFinalizada em: 23:15:13"

# Add new row 3 values
$ws.Range("A3").Value = "b"
$ws.Range("B3").Value = "b"
$ws.Range("C3").Value = "0:00:04.326423"
$ws.Range("D3").Value = "Finalizada"
$ws.Range("E3").Value = "Iniciada em: 23:16:06`nPausada em: 23:16:10`nFinalizada em: 23:16:16"
